# Updating to new dataset
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the stray header cell in A1 (it no longer belongs in the new dataset)
$ws.Range("A1").Clear()

# 2. Append a new organization row (19) for "Ingredion"
#    Copy the formatting from the last existing label cell (A18) so the new
#    label cell picks up the same bold/centered/bordered style.
$ws.Range("A18").Copy($ws.Range("A19"))
$ws.Range("A19").Value = "Ingredion"

# Fill the new row's data columns (B:AR) with 0 (no match yet for this org)
$ws.Range("B19:AR19").Value = 0
